$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the header row: columns were reorganized so that the
# "wingarea"/"span"/"mass" attributes move after the aero-coefficient
# columns (Ixx..cnr), and two new derivative columns (cldr/cndr) were
# introduced, replacing the old duplicate "cdr" header.
$ws.Range("A1").Value = "wingarea"
$ws.Range("B1").Value = "span"
$ws.Range("C1").Value = "mass"
$ws.Range("D1").Value = "Ixx"
$ws.Range("E1").Value = "Iyy"
$ws.Range("F1").Value = "Izz"
$ws.Range("G1").Value = "Ixz"
$ws.Range("H1").Value = "cyb"
$ws.Range("I1").Value = "cydr"
$ws.Range("J1").Value = "clb"
$ws.Range("K1").Value = "clp"
$ws.Range("L1").Value = "clr"
$ws.Range("M1").Value = "cldr"
$ws.Range("N1").Value = "cnb"
$ws.Range("O1").Value = "cnp"
$ws.Range("P1").Value = "cnr"
$ws.Range("Q1").Value = "cndr"

# Updated Ixx value
$ws.Range("D2").Value = 10970

# Leave the selection where the author left off
$ws.Range("D7").Select()
